# Updates coin price/volume data and fixes a row-order swap (Filecoin <-> FirstDigitalUSD)
# as published by the "Updated cryptos list" GitHub Actions job.
#
# Note: a leading single-quote is prefixed onto values that look like plain
# decimal numbers (e.g. "303.16", "1.00") so Excel stores them as text and
# keeps the exact original formatting instead of coercing them to a number.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "43.128.34"
$ws.Range("E2").Value = "  +1.65%  "
$ws.Range("D3").Value = "2.384.16"
$ws.Range("E3").Value = "  +4.27%  "
$ws.Range("E4").Value = "  +0.09%  "
$ws.Range("D5").Value = "'303.16"
$ws.Range("D6").Value = "'97.32"
$ws.Range("E6").Value = "  +2.03%  "
$ws.Range("E7").Value = "  +0.62%  "
$ws.Range("E8").Value = "  -0.05%  "
$ws.Range("D9").Value = "'0.502"
$ws.Range("E9").Value = "  +2.12%  "
$ws.Range("D10").Value = "'34.21"
$ws.Range("E10").Value = "  -0.44%  "
$ws.Range("D11").Value = "'0.0789"
$ws.Range("E11").Value = "  +0.46%  "
$ws.Range("E12").Value = "  +2.57%  "
$ws.Range("E13").Value = "  -1.97%  "
$ws.Range("D14").Value = "'6.79"
$ws.Range("E14").Value = "  +0.75%  "
$ws.Range("D15").Value = "2.759.60"
$ws.Range("E15").Value = "  +3.93%  "
$ws.Range("D16").Value = "2.404.09"
$ws.Range("E16").Value = "  +2.67%  "
$ws.Range("D17").Value = "'0.811"
$ws.Range("E17").Value = "  +4.26%  "
$ws.Range("D18").Value = "43.145.05"
$ws.Range("E18").Value = "  +1.81%  "
$ws.Range("E19").Value = "  +1.00%  "
$ws.Range("D20").Value = "'6.36"
$ws.Range("E20").Value = "  +6.54%  "
$ws.Range("D21").Value = "0.0₃0889"
$ws.Range("E21").Value = "  +0.29%  "
$ws.Range("D22").Value = "'68.43"
$ws.Range("E22").Value = "  +1.20%  "
$ws.Range("D23").Value = "'236.13"
$ws.Range("E23").Value = "  +0.28%  "
$ws.Range("D24").Value = "'2.22"
$ws.Range("E24").Value = "  -1.64%  "
$ws.Range("E25").Value = "  +1.92%  "
$ws.Range("E26").Value = "  -0.01%  "
$ws.Range("D27").Value = "'24.80"
$ws.Range("E27").Value = "  +2.24%  "
$ws.Range("E28").Value = "  +7.21%  "
$ws.Range("D29").Value = "'9.13"
$ws.Range("E29").Value = "  +1.25%  "
$ws.Range("D30").Value = "'31.57"
$ws.Range("E30").Value = "  -0.56%  "
$ws.Range("B31").Value = "FirstDigitalUSD"
$ws.Range("C31").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("D31").Value = "'1.00"
$ws.Range("E31").Value = "  +0.08%  "
$ws.Range("B32").Value = "Filecoin"
$ws.Range("C32").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D32").Value = "'5.11"
$ws.Range("E32").Value = "  +2.68%  "
$ws.Range("D33").Value = "'0.0729"
$ws.Range("E33").Value = "  +4.53%  "
$ws.Range("D34").Value = "'17.23"
$ws.Range("E34").Value = "  -1.15%  "
$ws.Range("D35").Value = "'1.86"
$ws.Range("E35").Value = "  +7.22%  "
$ws.Range("D36").Value = "'4.38"
$ws.Range("E36").Value = "  -0.60%  "
$ws.Range("D37").Value = "'2.30"
$ws.Range("E37").Value = "  -0.91%  "
$ws.Range("E38").Value = "  +2.29%  "
$ws.Range("D39").Value = "'2.81"
$ws.Range("E39").Value = "  +4.92%  "
$ws.Range("D40").Value = "'22.62"
$ws.Range("E40").Value = "  +13.12%  "
$ws.Range("E41").Value = "  +0.42%  "
$ws.Range("D42").Value = "'107.48"
$ws.Range("E42").Value = "  -35.07%  "
$ws.Range("D43").Value = "1.948.76"
$ws.Range("E43").Value = "  -0.79%  "
$ws.Range("D44").Value = "'0.0280"
$ws.Range("E44").Value = "  +1.05%  "
$ws.Range("E45").Value = "  +2.29%  "
$ws.Range("D46").Value = "'9.34"
$ws.Range("E46").Value = "  -10.54%  "
$ws.Range("D47").Value = "'2.75"
$ws.Range("E47").Value = "  +0.27%  "
$ws.Range("D48").Value = "2.611.09"
$ws.Range("E48").Value = "  +3.92%  "
$ws.Range("D49").Value = "'52.94"
$ws.Range("E49").Value = "  +0.08%  "
$ws.Range("D50").Value = "'72.29"
$ws.Range("E50").Value = "  +1.78%  "
$ws.Range("E51").Value = "  +1.89%  "
